$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header labels (legend text) in row 1
$ws.Range("B1").Value = "Potentially understandable by patients"
$ws.Range("C1").Value = "Not potentially understanably by patients"

# Row 1 now needs more vertical space for the longer wrapped text
$ws.Rows.Item(1).RowHeight = 82.5

# Move the active selection to C2
$ws.Range("C2").Select()
